$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.442.98'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.875.57'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7151'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07891'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3092'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08261'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.878.89'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7260'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.258'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.23'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').Value = '29.415.64'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.883'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '245.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007838'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.26'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '2.122.35'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('E22').Value = '  +6.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1612'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.016'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.30'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('E29').Value = '  -2.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.496'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.399'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.102'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05195'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.939'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range('E35').Value = '  +1.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7228'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.677'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01867'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('D40').Value = '1.184.71'
$ws.Range('E40').Value = '  +2.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9074'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.132'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.78'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.04'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5297'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '2.017.32'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.791'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.13%  '
$ws.Range('E49').Value = '  +6.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.290'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4304'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.09%  '
